# Updated cryptos list on Wed Mar 20 11:59:39 UTC 2024 with GitHub Actions
#
# Applies the latest scraped coinranking.com values to the "Price" (D) and
# "Volume(1h)" (E) columns of each crypto row, and reflects the rank swap
# between Toncoin and RenderToken (rows 23 & 24 change identity + data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = "63.377.31"
$ws.Cells.Item(2, 5).Value = "  +0.28%  "
# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = "3.285.71"
$ws.Cells.Item(3, 5).Value = "  +0.26%  "
# Row 4 - TetherUSD
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.15%  "
# Row 5 - BNB
$ws.Cells.Item(5, 4).Value = "533.04"
$ws.Cells.Item(5, 5).Value = "  +3.82%  "
# Row 6 - Solana
$ws.Cells.Item(6, 4).Value = "172.29"
$ws.Cells.Item(6, 5).Value = "  -3.10%  "
# Row 7 - XRP
$ws.Cells.Item(7, 4).Value = "0.596"
$ws.Cells.Item(7, 5).Value = "  +0.97%  "
# Row 8 - LidoStakedEther
$ws.Cells.Item(8, 4).Value = "3.287.10"
$ws.Cells.Item(8, 5).Value = "  +0.35%  "
# Row 9 - USDC (price unchanged)
$ws.Cells.Item(9, 5).Value = "  -0.16%  "
# Row 10 - Cardano
$ws.Cells.Item(10, 4).Value = "0.611"
$ws.Cells.Item(10, 5).Value = "  -1.06%  "
# Row 11 - Avalanche
$ws.Cells.Item(11, 4).Value = "53.51"
$ws.Cells.Item(11, 5).Value = "  -7.74%  "
# Row 12 - Dogecoin
$ws.Cells.Item(12, 4).Value = "0.135"
$ws.Cells.Item(12, 5).Value = "  +3.35%  "
# Row 13 - ShibaInu
$ws.Cells.Item(13, 4).Value = "0.0000258"
$ws.Cells.Item(13, 5).Value = "  +1.78%  "
# Row 14 - Polkadot
$ws.Cells.Item(14, 4).Value = "9.28"
$ws.Cells.Item(14, 5).Value = "  +1.77%  "
# Row 15 - WrappedliquidstakedEther2.0
$ws.Cells.Item(15, 4).Value = "3.793.92"
$ws.Cells.Item(15, 5).Value = "  +0.35%  "
# Row 16 - TRON
$ws.Cells.Item(16, 4).Value = "0.118"
$ws.Cells.Item(16, 5).Value = "  -1.21%  "
# Row 17 - WrappedEther
$ws.Cells.Item(17, 4).Value = "3.274.12"
$ws.Cells.Item(17, 5).Value = "  -0.19%  "
# Row 18 - Chainlink
$ws.Cells.Item(18, 4).Value = "17.43"
$ws.Cells.Item(18, 5).Value = "  +0.82%  "
# Row 19 - WrappedBTC
$ws.Cells.Item(19, 4).Value = "63.204.27"
$ws.Cells.Item(19, 5).Value = "  +0.11%  "
# Row 20 - Uniswap
$ws.Cells.Item(20, 4).Value = "11.19"
$ws.Cells.Item(20, 5).Value = "  +2.68%  "
# Row 21 - Polygon
$ws.Cells.Item(21, 4).Value = "0.965"
$ws.Cells.Item(21, 5).Value = "  +2.58%  "
# Row 22 - BitcoinCash
$ws.Cells.Item(22, 4).Value = "370.45"
$ws.Cells.Item(22, 5).Value = "  -0.33%  "

# Rows 23 & 24 swap identity: Toncoin overtakes RenderToken in rank, so the
# row that used to be Toncoin (23) now shows RenderToken's refreshed data,
# and the row that used to be RenderToken (24) now shows Toncoin's.
$ws.Cells.Item(23, 2).Value = "RenderToken"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(23, 4).Value = "11.48"
$ws.Cells.Item(23, 5).Value = "  +2.49%  "

$ws.Cells.Item(24, 2).Value = "Toncoin"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(24, 4).Value = "4.14"
$ws.Cells.Item(24, 5).Value = "  +10.20%  "

# Row 25 - PancakeSwap
$ws.Cells.Item(25, 4).Value = "3.77"
$ws.Cells.Item(25, 5).Value = "  +3.67%  "
# Row 26 - Litecoin
$ws.Cells.Item(26, 4).Value = "81.28"
$ws.Cells.Item(26, 5).Value = "  +2.27%  "
# Row 27 - LEO (price unchanged)
$ws.Cells.Item(27, 5).Value = "  +4.12%  "
# Row 28 - ImmutableX
$ws.Cells.Item(28, 4).Value = "2.67"
$ws.Cells.Item(28, 5).Value = "  +1.56%  "
# Row 29 - InternetComputer(DFINITY)
$ws.Cells.Item(29, 4).Value = "11.33"
$ws.Cells.Item(29, 5).Value = "  +0.39%  "
# Row 30 - Filecoin
$ws.Cells.Item(30, 4).Value = "8.28"
$ws.Cells.Item(30, 5).Value = "  +0.03%  "
# Row 31 - EthereumClassic
$ws.Cells.Item(31, 4).Value = "28.73"
$ws.Cells.Item(31, 5).Value = "  +1.70%  "
# Row 32 - Bittensor (force text format so the trailing zero in "642.70"
# survives - plain Value assignment would coerce it to the number 642.7)
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "642.70"
$ws.Cells.Item(32, 5).Value = "  -0.58%  "
# Row 33 - NEARProtocol
$ws.Cells.Item(33, 4).Value = "6.45"
$ws.Cells.Item(33, 5).Value = "  -3.51%  "
# Row 34 - Cosmos
$ws.Cells.Item(34, 4).Value = "11.28"
$ws.Cells.Item(34, 5).Value = "  +1.28%  "
# Row 35 - Hedera (price unchanged)
$ws.Cells.Item(35, 5).Value = "  +3.47%  "
# Row 36 - OKB
$ws.Cells.Item(36, 4).Value = "56.92"
$ws.Cells.Item(36, 5).Value = "  -2.59%  "
# Row 37 - Dai (price unchanged)
$ws.Cells.Item(37, 5).Value = "  +0.13%  "
# Row 38 - InjectiveProtocol
$ws.Cells.Item(38, 4).Value = "36.79"
$ws.Cells.Item(38, 5).Value = "  +2.32%  "
# Row 39 - TheGraph (price unchanged)
$ws.Cells.Item(39, 5).Value = "  +0.35%  "
# Row 40 - PEPE
$ws.Cells.Item(40, 4).Value = "0.0₃0737"
$ws.Cells.Item(40, 5).Value = "  +11.57%  "
# Row 41 - FirstDigitalUSD
$ws.Cells.Item(41, 4).Value = "0.996"
$ws.Cells.Item(41, 5).Value = "  -0.17%  "
# Row 42 - Fetch.AI
$ws.Cells.Item(42, 4).Value = "2.63"
$ws.Cells.Item(42, 5).Value = "  +7.35%  "
# Row 43 - Kaspa (price unchanged)
$ws.Cells.Item(43, 5).Value = "  -1.65%  "
# Row 44 - Maker
$ws.Cells.Item(44, 4).Value = "2.908.05"
$ws.Cells.Item(44, 5).Value = "  +0.88%  "
# Row 45 - Stacks
$ws.Cells.Item(45, 4).Value = "2.93"
$ws.Cells.Item(45, 5).Value = "  +2.61%  "
# Row 46 - WEMIXToken
$ws.Cells.Item(46, 4).Value = "2.68"
$ws.Cells.Item(46, 5).Value = "  +4.79%  "
# Row 47 - VeChain
$ws.Cells.Item(47, 4).Value = "0.0399"
$ws.Cells.Item(47, 5).Value = "  +4.03%  "
# Row 48 - ThetaToken (price unchanged)
$ws.Cells.Item(48, 5).Value = "  -1.79%  "
# Row 49 - ApeXProtocol
$ws.Cells.Item(49, 4).Value = "3.04"
$ws.Cells.Item(49, 5).Value = "  +4.36%  "
# Row 50 - Stellar
$ws.Cells.Item(50, 4).Value = "0.126"
$ws.Cells.Item(50, 5).Value = "  +2.46%  "
# Row 51 - Monero
$ws.Cells.Item(51, 4).Value = "134.54"
$ws.Cells.Item(51, 5).Value = "  +3.39%  "
